# Add a new "2022-Q4" sheet right after "总计" (i.e. before the current
# "2022-Q3" sheet), shift the existing quarter sheets down, and add the
# corresponding summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)          # "总计"
$q3Old = $wb.Worksheets.Item(2)           # current "2022-Q3" (will stay in place, content unchanged)

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet before the current "2022-Q3" tab
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($q3Old)

# Copy the header-row style (bold + border, style index used by B1:H1 on
# every quarter sheet) from the existing "2022-Q3" sheet before we put any
# data in, then fill in the Q4 header text.
$q3Old.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Row data: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0,  "501208", "中欧创新未来混合（LOF）",              "54.76", "85.24", "2.96", "1.6209", 10),
    @(1,  "007484", "信澳核心科技混合",                      "21.50", "93.51", "5.58", "1.1997", 1),
    @(2,  "005763", "中欧电子信息产业沪港深股票C",          "14.81", "91.56", "5.88", "0.8708", 3),
    @(3,  "001513", "易方达信息产业混合",                    "33.11", "90.45", "1.93", "0.6390", 7),
    @(4,  "010013", "易方达信息行业精选股票",                "22.89", "88.67", "2.13", "0.4876", 8),
    @(5,  "004616", "中欧电子信息产业沪港深股票A",          "6.80",  "91.56", "5.88", "0.3998", 3),
    @(6,  "506002", "易方达科创板两年定期开放混合",          "14.41", "93.98", "2.68", "0.3862", 10),
    @(7,  "010622", "恒越成长精选混合A",                    "10.54", "68.17", "2.14", "0.2256", 6),
    @(8,  "001411", "诺安创新驱动灵活配置混合A",            "7.22",  "91.15", "2.28", "0.1646", 6),
    @(9,  "010824", "天弘创新成长混合A",                    "2.73",  "82.47", "5.35", "0.1461", 3),
    @(10, "002051", "诺安创新驱动灵活配置混合C",            "4.60",  "91.15", "2.28", "0.1049", 6),
    @(11, "010623", "恒越成长精选混合C",                    "4.46",  "68.17", "2.14", "0.0954", 6),
    @(12, "010825", "天弘创新成长混合C",                    "1.01",  "82.47", "5.35", "0.0540", 3),
    @(13, "015919", "申万菱信专精特新主题混合A",            "0.39",  "48.55", "4.83", "0.0188", 2),
    @(14, "011214", "招商惠润一年定期开放混合（MOM）A",      "0.48",  "68.20", "1.70", "0.0082", 10),
    @(15, "016380", "华宝专精特新混合A",                    "0.10",  "90.16", "4.06", "0.0041", 7),
    @(16, "015920", "申万菱信专精特新主题混合C",            "0.02",  "48.55", "4.83", "0.0010", 2),
    @(17, "011215", "招商惠润一年定期开放混合（MOM）C",      "0.06",  "68.20", "1.70", "0.0010", 10),
    @(18, "016381", "华宝专精特新混合C",                    "0.01",  "90.16", "4.06", "0.0004", 7)
)

foreach ($r in $rows) {
    $rowNum = $r[0] + 2

    # Column A: numeric row index, styled like the other quarter sheets'
    # index column (bold/centered/bordered).
    $q3Old.Cells.Item(2,1).Copy()
    $q4.Cells.Item($rowNum,1).PasteSpecial(-4122)
    $q4.Cells.Item($rowNum,1).Value = $r[0]

    # Column B: fund code -- keep as text even though it looks numeric.
    $q4.Cells.Item($rowNum,2).NumberFormat = "@"
    $q4.Cells.Item($rowNum,2).Value = $r[1]

    # Column C: fund name (plain text).
    $q4.Cells.Item($rowNum,3).Value = $r[2]

    # Columns D,E,F,G: decimal-looking numbers stored as text.
    $q4.Cells.Item($rowNum,4).NumberFormat = "@"
    $q4.Cells.Item($rowNum,4).Value = $r[3]

    $q4.Cells.Item($rowNum,5).NumberFormat = "@"
    $q4.Cells.Item($rowNum,5).Value = $r[4]

    $q4.Cells.Item($rowNum,6).NumberFormat = "@"
    $q4.Cells.Item($rowNum,6).Value = $r[5]

    $q4.Cells.Item($rowNum,7).NumberFormat = "@"
    $q4.Cells.Item($rowNum,7).Value = $r[6]

    # Column H: position rank, numeric.
    $q4.Cells.Item($rowNum,8).Value = $r[7]
}

# Name the sheet last -- renaming before the cross-sheet copy/paste above
# confused style application in this runtime.
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet. The data in columns B/C/D cascades
#    down by one row (each row now shows the previous row's quarter), a
#    brand-new 2022-Q4 row is written into row 2, and one new row is
#    appended at the bottom carrying forward the old last row's data.
#    Column A (the plain row counter) is left untouched for rows 2-5 and
#    simply continues (+1) for the new row 6.
# ---------------------------------------------------------------------

# Capture the existing rows 2-5 (B/C/D) before overwriting anything.
$oldRows = @()
foreach ($r in 2..5) {
    $oldRows += ,@($total.Cells.Item($r,2).Value, $total.Cells.Item($r,3).Value, $total.Cells.Item($r,4).Value)
}
$lastIndex = $total.Cells.Item(5,1).Value

# New row 6: same style as the existing index rows, carrying the old row
# 5 (2021-Q4) data forward, with the counter continuing on from row 5.
$total.Cells.Item(5,1).Copy()
$total.Cells.Item(6,1).PasteSpecial(-4122)
$total.Cells.Item(6,1).Value = $lastIndex + 1
$total.Cells.Item(6,2).Value = $oldRows[3][0]
$total.Cells.Item(6,3).Value = $oldRows[3][1]
$total.Cells.Item(6,4).Value = $oldRows[3][2]

# Rows 5,4,3 take on the data that used to belong to rows 4,3,2.
$total.Cells.Item(5,2).Value = $oldRows[2][0]
$total.Cells.Item(5,3).Value = $oldRows[2][1]
$total.Cells.Item(5,4).Value = $oldRows[2][2]

$total.Cells.Item(4,2).Value = $oldRows[1][0]
$total.Cells.Item(4,3).Value = $oldRows[1][1]
$total.Cells.Item(4,4).Value = $oldRows[1][2]

$total.Cells.Item(3,2).Value = $oldRows[0][0]
$total.Cells.Item(3,3).Value = $oldRows[0][1]
$total.Cells.Item(3,4).Value = $oldRows[0][2]

# Row 2 gets the brand-new 2022-Q4 figures.
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 19
$total.Cells.Item(2,4).Value = 6.43
